$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: propagate number/date formatting from column F (the old column D) into new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 3: populate the new D and E columns with the new period data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4344200
$ws.Range("E8").Value = 4770800
$ws.Range("D9").Value = 3537500
$ws.Range("E9").Value = 3778300
$ws.Range("D10").Value = 806700
$ws.Range("E10").Value = 992500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 113700
$ws.Range("E15").Value = 108900
$ws.Range("D17").Value = 4126600
$ws.Range("E17").Value = 4267300
$ws.Range("D18").Value = 217600
$ws.Range("E18").Value = 503500
$ws.Range("D20").Value = 8100
$ws.Range("E20").Value = 7300
$ws.Range("D21").Value = 339500
$ws.Range("E21").Value = 619700
$ws.Range("D22").Value = 33900
$ws.Range("E22").Value = 32400
$ws.Range("D23").Value = 191800
$ws.Range("E23").Value = 478400
$ws.Range("D24").Value = 36300
$ws.Range("E24").Value = 116300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 155500
$ws.Range("E26").Value = 362100
$ws.Range("D27").Value = 133700
$ws.Range("E27").Value = 341300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 7800
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -8100
$ws.Range("E32").Value = -7300
$ws.Range("D33").Value = 141500
$ws.Range("E33").Value = 341300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 141500
$ws.Range("E35").Value = 341300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1154800
$ws.Range("E41").Value = 1075700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 705700
$ws.Range("E43").Value = 862500
$ws.Range("D44").Value = 1354400
$ws.Range("E44").Value = 1862500
$ws.Range("D45").Value = 81500
$ws.Range("E45").Value = 57100
$ws.Range("D46").Value = 3296400
$ws.Range("E46").Value = 3857700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 4682500
$ws.Range("E48").Value = 4673400
$ws.Range("D49").Value = 2673500
$ws.Range("E49").Value = 2689800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 342200
$ws.Range("E52").Value = 274200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 10994600
$ws.Range("E54").Value = 11495100
$ws.Range("D57").Value = 872600
$ws.Range("E57").Value = 1176400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 295500
$ws.Range("E59").Value = 364200
$ws.Range("D60").Value = 1168200
$ws.Range("E60").Value = 1540600
$ws.Range("D61").Value = 2411500
$ws.Range("E61").Value = 2409100
$ws.Range("D62").Value = 955800
$ws.Range("E62").Value = 981800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5076000
$ws.Range("E66").Value = 5482800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 4196900
$ws.Range("E72").Value = 4112600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 5918600
$ws.Range("E76").Value = 6012300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 141500
$ws.Range("E81").Value = 341300
$ws.Range("D83").Value = 113700
$ws.Range("E83").Value = 108900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 424500
$ws.Range("E89").Value = 401800
$ws.Range("D91").Value = -13000
$ws.Range("E91").Value = -9500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -101400
$ws.Range("E94").Value = -113100
$ws.Range("D96").Value = -57600
$ws.Range("E96").Value = -58400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -240900
$ws.Range("E100").Value = -193700
$ws.Range("D101").Value = -3200
$ws.Range("E101").Value = 800
$ws.Range("D102").Value = 79100
$ws.Range("E102").Value = 95800

# Step 4: apply corrected values (restated figures) on top of the shifted historical columns
$ws.Range("H9").Value = 3107100
$ws.Range("H10").Value = 885600
$ws.Range("H17").Value = 3594200
$ws.Range("H18").Value = 398500
$ws.Range("H20").Value = 4800
$ws.Range("H32").Value = -4800
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("F91").Value = -19000
$ws.Range("G91").Value = -12600
$ws.Range("H91").Value = -14100
$ws.Range("I91").Value = -10200
$ws.Range("J91").Value = -20500
